# edit.ps1 - apply the "Added assortativity to info 2012" commit to the deck
#
# Content changes made by this script:
#  1. Slide 2 ("Outline"), bullet "Implement k-dense in igraph" -> capitalize
#     the library name so it reads "...in IGraph" (the run gets split into
#     "IG" + "raph" exactly like PowerPoint's autocorrect/retype would do).
#  2. Slide 7 ("Implementing k-dense") title placeholder loses its now
#     redundant trailing <a:endParaRPr/> - reproduced by clearing and
#     retyping the title text so the paragraph is rebuilt without it.

$p = $ppt.ActivePresentation

# --- 1. "igraph" -> "IGraph" on the Outline slide -------------------------
$outlineSlide = $p.Slides.Item(2)
$outlineShape = $outlineSlide.Shapes.Item(2)
$outlineRange = $outlineShape.TextFrame2.TextRange

$fullText = $outlineRange.Text
$matchIndex = $fullText.IndexOf("igraph")
if ($matchIndex -ge 0) {
    $startPos = $matchIndex + 1   # TextRange2.Characters is 1-based
    $prefix = $outlineRange.Characters($startPos, 2)
    $prefix.Text = "IG"
}

# --- 2. Drop the stray endParaRPr on the "Implementing k-dense" title -----
$implSlide = $p.Slides.Item(7)
$titleShape = $implSlide.Shapes.Item(1)
$titleRange = $titleShape.TextFrame2.TextRange

$titleText = $titleRange.Text
$titleRange.Delete()
$titleRange.Text = $titleText
